# Tutorial 1 - Introduction.pptx
# Add "Andrew Tan" (the tutor's name) plus a trailing blank line to the
# subtitle placeholder of the title slide (Slide 1, Shape 2 - "Subtitle 2"),
# which currently just reads "AY18/19 - week 3".
#
# Resulting paragraphs inside the subtitle text box:
#   1) AY18/19 - week 3      (existing - untouched)
#   2) Andrew Tan            (new)
#   3) <empty paragraph>     (new)

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# Append a new paragraph with the tutor's name, followed by one more
# (empty) paragraph break, to the end of the existing text.
$tr.InsertAfter([char]13 + "Andrew Tan" + [char]13)
